# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): labels with the same bold/centered/bordered style
# used by the existing header cells (e.g. A1 uses style index 1). ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting from an existing styled header cell (A1) onto
# the three new header cells, so they pick up the bold/border/center style.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-60): season record repeated for every player/row. ---
$ws.Range("AD2:AD60").Value = 90
$ws.Range("AE2:AE60").Value = 72
$ws.Range("AF2:AF60").Value = 0
